{"js": "// Replace the worksheet date and each two-digit-by-two-digit multiplication\n// equation with its updated value, per the commit's regenerated answer key.\n// Every \"before\" string below is unique within the document, so an exact,\n// case-sensitive literal search safely targets the single matching run.\nconst replacements = [\n  [\"2024-12-23 Monday\", \"2024-12-24 Tuesday\"],\n  [\"41\u00d734=1394\", \"90\u00d717=1530\"],\n  [\"71\u00d751=3621\", \"58\u00d784=4872\"],\n  [\"11\u00d749=539\", \"56\u00d779=4424\"],\n  [\"14\u00d756=784\", \"59\u00d758=3422\"],\n  [\"63\u00d729=1827\", \"33\u00d785=2805\"],\n  [\"13\u00d721=273\", \"75\u00d796=7200\"],\n  [\"24\u00d738=912\", \"81\u00d727=2187\"],\n  [\"82\u00d724=1968\", \"22\u00d718=396\"],\n  [\"84\u00d792=7728\", \"57\u00d791=5187\"],\n  [\"94\u00d786=8084\", \"36\u00d779=2844\"],\n  [\"40\u00d718=720\", \"19\u00d750=950\"],\n  [\"43\u00d748=2064\", \"27\u00d722=594\"],\n  [\"40\u00d729=1160\", \"18\u00d790=1620\"],\n  [\"73\u00d724=1752\", \"63\u00d758=3654\"],\n  [\"30\u00d749=1470\", \"41\u00d712=492\"],\n  [\"21\u00d719=399\", \"96\u00d771=6816\"],\n  [\"79\u00d762=4898\", \"70\u00d766=4620\"],\n  [\"42\u00d716=672\", \"99\u00d724=2376\"],\n  [\"94\u00d779=7426\", \"49\u00d795=4655\"],\n  [\"81\u00d782=6642\", \"19\u00d742=798\"],\n  [\"29\u00d758=1682\", \"40\u00d755=2200\"],\n  [\"54\u00d791=4914\", \"44\u00d716=704\"],\n  [\"61\u00d762=3782\", \"25\u00d769=1725\"],\n  [\"36\u00d711=396\", \"20\u00d782=1640\"],\n  [\"95\u00d792=8740\", \"36\u00d799=3564\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2024-12-23 Monday\", \"2024-12-24 Tuesday\"),\n  @(\"41\u00d734=1394\", \"90\u00d717=1530\"),\n  @(\"71\u00d751=3621\", \"58\u00d784=4872\"),\n  @(\"11\u00d749=539\", \"56\u00d779=4424\"),\n  @(\"14\u00d756=784\", \"59\u00d758=3422\"),\n  @(\"63\u00d729=1827\", \"33\u00d785=2805\"),\n  @(\"13\u00d721=273\", \"75\u00d796=7200\"),\n  @(\"24\u00d738=912\", \"81\u00d727=2187\"),\n  @(\"82\u00d724=1968\", \"22\u00d718=396\"),\n  @(\"84\u00d792=7728\", \"57\u00d791=5187\"),\n  @(\"94\u00d786=8084\", \"36\u00d779=2844\"),\n  @(\"40\u00d718=720\", \"19\u00d750=950\"),\n  @(\"43\u00d748=2064\", \"27\u00d722=594\"),\n  @(\"40\u00d729=1160\", \"18\u00d790=1620\"),\n  @(\"73\u00d724=1752\", \"63\u00d758=3654\"),\n  @(\"30\u00d749=1470\", \"41\u00d712=492\"),\n  @(\"21\u00d719=399\", \"96\u00d771=6816\"),\n  @(\"79\u00d762=4898\", \"70\u00d766=4620\"),\n  @(\"42\u00d716=672\", \"99\u00d724=2376\"),\n  @(\"94\u00d779=7426\", \"49\u00d795=4655\"),\n  @(\"81\u00d782=6642\", \"19\u00d742=798\"),\n  @(\"29\u00d758=1682\", \"40\u00d755=2200\"),\n  @(\"54\u00d791=4914\", \"44\u00d716=704\"),\n  @(\"61\u00d762=3782\", \"25\u00d769=1725\"),\n  @(\"36\u00d711=396\", \"20\u00d782=1640\"),\n  @(\"95\u00d792=8740\", \"36\u00d799=3564\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n  #         MatchAllWordForms, Forward, Wrap:=wdFindContinue(1), Format, ReplaceWith,\n  #         Replace:=wdReplaceAll(2)) \u2014 each \"before\" string is unique in the\n  # document, so a single case-sensitive literal Find/Replace-all call targets\n  # exactly the one matching run.\n  $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    throw \"Could not find text to replace: $oldText\"\n  }\n}\n"}
